$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 645.25
$ws.Range("I2").Value = 700.3333
$ws.Range("J2").Value = 480
$ws.Range("K2").Value = 700.3333
$ws.Range("L2").Value = 480
$ws.Range("M2").Value = -587.3333
$ws.Range("N2").Value = -706
$ws.Range("H21").Value = 10884.538
$ws.Range("I21").Value = 7687.5
$ws.Range("J21").Value = 15999.8
$ws.Range("K21").Value = 7687.5
$ws.Range("L21").Value = 15999.8
$ws.Range("M21").Value = -7219.5
$ws.Range("N21").Value = -16935.8
$ws.Range("H23").Value = 10884.538
$ws.Range("I23").Value = 7687.5
$ws.Range("J23").Value = 15999.8
$ws.Range("K23").Value = 7687.5
$ws.Range("L23").Value = 15999.8
$ws.Range("M23").Value = -7453.5
$ws.Range("N23").Value = -16467.8
$ws.Range("H29").Value = 915.4286
$ws.Range("J29").Value = 704
$ws.Range("L29").Value = 2112
$ws.Range("N29").Value = -2674
$ws.Range("H38").Value = 3135.3684
$ws.Range("I38").Value = 160.25
$ws.Range("J38").Value = 5299.091
$ws.Range("K38").Value = 480.75
$ws.Range("L38").Value = 15897.273
$ws.Range("M38").Value = -108.75
$ws.Range("N38").Value = -16641.273
$ws.Range("H58").Value = 10759.4
$ws.Range("I58").Value = 1398
$ws.Range("J58").Value = 14771.429
$ws.Range("K58").Value = 4194
$ws.Range("L58").Value = 44314.287
$ws.Range("M58").Value = -4044
$ws.Range("N58").Value = -44614.287
$ws.Range("H87").Value = 22042
$ws.Range("J87").Value = 22042
$ws.Range("L87").Value = 22042
$ws.Range("N87").Value = -24538
$ws.Range("H90").Value = 22042
$ws.Range("J90").Value = 22042
$ws.Range("L90").Value = 66126
$ws.Range("N90").Value = -78606
$ws.Range("H93").Value = 22711.111
$ws.Range("J93").Value = 22711.111
$ws.Range("L93").Value = 22711.111
$ws.Range("N93").Value = -27703.111
$ws.Range("H98").Value = 5842.0566
$ws.Range("I98").Value = 4623.577
$ws.Range("J98").Value = 7015.407
$ws.Range("K98").Value = 4623.577
$ws.Range("L98").Value = 7015.407
$ws.Range("M98").Value = -3125.577
$ws.Range("N98").Value = -10011.407
$ws.Range("H122").Value = 5842.0566
$ws.Range("I122").Value = 4623.577
$ws.Range("J122").Value = 7015.407
$ws.Range("K122").Value = 13870.731
$ws.Range("L122").Value = 21046.221
$ws.Range("M122").Value = -11420.731
$ws.Range("N122").Value = -25946.221
$ws.Range("H127").Value = 975.3333
$ws.Range("I127").Value = 797
$ws.Range("J127").Value = 1153.6666
$ws.Range("K127").Value = 2391
$ws.Range("L127").Value = 3460.9998
$ws.Range("M127").Value = 2569
$ws.Range("N127").Value = -13380.9998
$ws.Range("H131").Value = 2792.5
$ws.Range("I131").Value = 2046.6666
$ws.Range("J131").Value = 4284.1665
$ws.Range("K131").Value = 6139.9998
$ws.Range("L131").Value = 12852.4995
$ws.Range("M131").Value = -1099.9998
$ws.Range("N131").Value = -22932.4995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1250
$ws.Range("I102").Value = 1000
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 1000
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = 622
$ws.Range("N102").Value = -4744
$ws.Range("H103").Value = 34723.08
$ws.Range("J103").Value = 34723.08
$ws.Range("L103").Value = 34723.08
$ws.Range("N103").Value = -37067.08
$ws.Range("H122").Value = 3368.7
$ws.Range("I122").Value = 2886.7026
$ws.Range("J122").Value = 4740.5386
$ws.Range("K122").Value = 8660.1078
$ws.Range("L122").Value = 14221.6158
$ws.Range("M122").Value = -6210.1078
$ws.Range("N122").Value = -19121.6158

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 33000
$ws.Range("J95").Value = 33000
$ws.Range("L95").Value = 33000
$ws.Range("N95").Value = -38492
$ws.Range("H103").Value = 36653.848
$ws.Range("J103").Value = 36653.848
$ws.Range("L103").Value = 36653.848
$ws.Range("N103").Value = -38997.848

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H98").Value = 47107.5
$ws.Range("J98").Value = 47107.5
$ws.Range("L98").Value = 47107.5
$ws.Range("N98").Value = -51599.5
$ws.Range("H122").Value = 2835.0833
$ws.Range("I122").Value = 1424.2
$ws.Range("J122").Value = 3842.8572
$ws.Range("K122").Value = 4272.6
$ws.Range("L122").Value = 11528.5716
$ws.Range("M122").Value = -1822.6
$ws.Range("N122").Value = -16428.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 5435354.5
$ws.Range("I113").Value = 650.2222
$ws.Range("J113").Value = 8929093
$ws.Range("K113").Value = 1950.6666
$ws.Range("L113").Value = 26787279
$ws.Range("M113").Value = 219.3334
$ws.Range("N113").Value = -26791619
$ws.Range("H131").Value = 771.6061
$ws.Range("J131").Value = 818.2222
$ws.Range("L131").Value = 2454.6666
$ws.Range("N131").Value = -12534.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 11770554
$ws.Range("I11").Value = 27000000
$ws.Range("J11").Value = 2252150.5
$ws.Range("K11").Value = 27000000
$ws.Range("L11").Value = 2252150.5
$ws.Range("M11").Value = -26999861
$ws.Range("N11").Value = -2252428.5
$ws.Range("H102").Value = 2533.6
$ws.Range("I102").Value = 1691.4667
$ws.Range("K102").Value = 1691.4667
$ws.Range("M102").Value = -69.46669999999995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5453.385
$ws.Range("I7").Value = 2268.8
$ws.Range("J7").Value = 7443.75
$ws.Range("K7").Value = 2268.8
$ws.Range("L7").Value = 7443.75
$ws.Range("M7").Value = -2156.8
$ws.Range("N7").Value = -7667.75
$ws.Range("H22").Value = 60703.94
$ws.Range("I22").Value = 101446.9
$ws.Range("J22").Value = 2499.7144
$ws.Range("K22").Value = 101446.9
$ws.Range("L22").Value = 2499.7144
$ws.Range("M22").Value = -101151.9
$ws.Range("N22").Value = -3089.7144
$ws.Range("H27").Value = 60703.94
$ws.Range("I27").Value = 101446.9
$ws.Range("J27").Value = 2499.7144
$ws.Range("K27").Value = 101446.9
$ws.Range("L27").Value = 2499.7144
$ws.Range("M27").Value = -101339.9
$ws.Range("N27").Value = -2713.7144
$ws.Range("H40").Value = 10840.9
$ws.Range("I40").Value = 9959.799999999999
$ws.Range("J40").Value = 11722
$ws.Range("K40").Value = 9959.799999999999
$ws.Range("L40").Value = 11722
$ws.Range("M40").Value = -9823.799999999999
$ws.Range("N40").Value = -11994
$ws.Range("H126").Value = 5453.385
$ws.Range("I126").Value = 2268.8
$ws.Range("J126").Value = 7443.75
$ws.Range("K126").Value = 6806.400000000001
$ws.Range("L126").Value = 22331.25
$ws.Range("M126").Value = -4336.400000000001
$ws.Range("N126").Value = -27271.25
$ws.Range("H135").Value = 50666.668
$ws.Range("J135").Value = 50666.668
$ws.Range("L135").Value = 50666.668
$ws.Range("N135").Value = -60806.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 39700
$ws.Range("J80").Value = 39700
$ws.Range("L80").Value = 39700
$ws.Range("N80").Value = -41696
$ws.Range("H81").Value = 2469.1538
$ws.Range("I81").Value = 2566.5833
$ws.Range("J81").Value = 1300
$ws.Range("K81").Value = 5133.1666
$ws.Range("L81").Value = 2600
$ws.Range("M81").Value = -4072.1666
$ws.Range("N81").Value = -4722
$ws.Range("H83").Value = 39700
$ws.Range("J83").Value = 39700
$ws.Range("L83").Value = 119100
$ws.Range("N83").Value = -129084
$ws.Range("H84").Value = 2469.1538
$ws.Range("I84").Value = 2566.5833
$ws.Range("J84").Value = 1300
$ws.Range("K84").Value = 25665.833
$ws.Range("L84").Value = 13000
$ws.Range("M84").Value = -20361.833
$ws.Range("N84").Value = -23608
$ws.Range("H101").Value = 11460.2
$ws.Range("J101").Value = 11460.2
$ws.Range("L101").Value = 11460.2
$ws.Range("N101").Value = -17950.2
$ws.Range("H103").Value = 32798
$ws.Range("J103").Value = 32798
$ws.Range("L103").Value = 32798
$ws.Range("N103").Value = -35142
$ws.Range("H113").Value = 204.9
$ws.Range("I113").Value = 206.61111
$ws.Range("J113").Value = 189.5
$ws.Range("K113").Value = 619.8333299999999
$ws.Range("L113").Value = 568.5
$ws.Range("M113").Value = 1550.16667
$ws.Range("N113").Value = -4908.5
